$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item(2, 4) "304.27"
Set-TextValue $ws.Cells.Item(2, 5) "4.01%"
Set-TextValue $ws.Cells.Item(3, 4) "35.73"
Set-TextValue $ws.Cells.Item(3, 5) "14.46%"
Set-TextValue $ws.Cells.Item(4, 4) "5.064"
Set-TextValue $ws.Cells.Item(4, 5) "1.88%"
Set-TextValue $ws.Cells.Item(5, 4) "0.07830"
Set-TextValue $ws.Cells.Item(5, 5) "4.70%"
Set-TextValue $ws.Cells.Item(6, 4) "2.259"
Set-TextValue $ws.Cells.Item(6, 5) "-2.12%"
Set-TextValue $ws.Cells.Item(7, 4) "8.118"
Set-TextValue $ws.Cells.Item(7, 5) "4.45%"
Set-TextValue $ws.Cells.Item(8, 4) "4.003"
Set-TextValue $ws.Cells.Item(8, 5) "6.03%"
Set-TextValue $ws.Cells.Item(9, 4) "0.9279"
Set-TextValue $ws.Cells.Item(9, 5) "0.96%"
Set-TextValue $ws.Cells.Item(10, 4) "0.09818"
Set-TextValue $ws.Cells.Item(10, 5) "3.89%"
Set-TextValue $ws.Cells.Item(11, 4) "0.1821"
Set-TextValue $ws.Cells.Item(11, 5) "5.12%"
Set-TextValue $ws.Cells.Item(12, 4) "0.08662"
Set-TextValue $ws.Cells.Item(12, 5) "3.78%"
Set-TextValue $ws.Cells.Item(13, 4) "0.03410"
Set-TextValue $ws.Cells.Item(13, 5) "3.54%"
Set-TextValue $ws.Cells.Item(14, 4) "0.09926"
Set-TextValue $ws.Cells.Item(14, 5) "-0.14%"
Set-TextValue $ws.Cells.Item(15, 4) "0.001496"
Set-TextValue $ws.Cells.Item(15, 5) "0.02%"
Set-TextValue $ws.Cells.Item(16, 4) "0.005729"
Set-TextValue $ws.Cells.Item(16, 5) "0.21%"
Set-TextValue $ws.Cells.Item(17, 4) "3.486"
Set-TextValue $ws.Cells.Item(17, 5) "0.41%"
Set-TextValue $ws.Cells.Item(18, 5) "-3.11%"
Set-TextValue $ws.Cells.Item(19, 5) "2.99%"
Set-TextValue $ws.Cells.Item(20, 5) "0.90%"
Set-TextValue $ws.Cells.Item(21, 4) "4.548"
Set-TextValue $ws.Cells.Item(21, 5) "11.23%"
Set-TextValue $ws.Cells.Item(22, 5) "5.33%"
Set-TextValue $ws.Cells.Item(23, 4) "0.04669"
Set-TextValue $ws.Cells.Item(23, 5) "3.25%"
Set-TextValue $ws.Cells.Item(24, 4) "0.001238"
Set-TextValue $ws.Cells.Item(24, 5) "1.44%"
Set-TextValue $ws.Cells.Item(25, 4) "0.004502"
Set-TextValue $ws.Cells.Item(25, 5) "4.55%"
Set-TextValue $ws.Cells.Item(26, 4) "0.0001299"
Set-TextValue $ws.Cells.Item(26, 5) "0.01%"
Set-TextValue $ws.Cells.Item(27, 4) "0.0002698"
Set-TextValue $ws.Cells.Item(27, 5) "-20.48%"
Set-TextValue $ws.Cells.Item(39, 4) "0.01762"
Set-TextValue $ws.Cells.Item(39, 5) "8.79%"
Set-TextValue $ws.Cells.Item(40, 4) "0.04704"
Set-TextValue $ws.Cells.Item(40, 5) "2.96%"
Set-TextValue $ws.Cells.Item(41, 4) "0.007975"
Set-TextValue $ws.Cells.Item(41, 5) "6.30%"
Set-TextValue $ws.Cells.Item(42, 4) "0.1420"
Set-TextValue $ws.Cells.Item(42, 5) "4.38%"
Set-TextValue $ws.Cells.Item(43, 4) "0.008138"
Set-TextValue $ws.Cells.Item(43, 5) "-17.23%"
Set-TextValue $ws.Cells.Item(44, 4) "0.002299"
Set-TextValue $ws.Cells.Item(44, 5) "6.56%"
Set-TextValue $ws.Cells.Item(45, 4) "0.009121"
Set-TextValue $ws.Cells.Item(45, 5) "0.90%"
Set-TextValue $ws.Cells.Item(46, 4) "0.00006153"
Set-TextValue $ws.Cells.Item(46, 5) "0.85%"
Set-TextValue $ws.Cells.Item(47, 5) "0.01%"
Set-TextValue $ws.Cells.Item(48, 4) "4.051"
Set-TextValue $ws.Cells.Item(48, 5) "44.97%"
Set-TextValue $ws.Cells.Item(49, 4) "0.002689"
Set-TextValue $ws.Cells.Item(49, 5) "34.51%"
Set-TextValue $ws.Cells.Item(50, 4) "0.00002099"
Set-TextValue $ws.Cells.Item(50, 5) "0.01%"
Set-TextValue $ws.Cells.Item(51, 4) "0.0001999"
Set-TextValue $ws.Cells.Item(51, 5) "0.01%"
